# Insert a new data row at row 497 (pushing existing rows 497:561 down to 498:562)
# and populate it with the new price-report entry for "Región de Ñuble".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("497:497").Insert()

$ws.Range("A497").Value = 5
$ws.Range("B497").Value = "Macroferia Regional de Talca"
$ws.Range("C497").Value = "Maule"
$ws.Range("D497").Value = 45124
$ws.Range("E497").Value = 7
$ws.Range("F497").Value = 100114013
$ws.Range("G497").Value = "Zanahoria"
$ws.Range("H497").Value = "Sin especificar"
$ws.Range("I497").Value = "Primera"
$ws.Range("J497").Value = 500
$ws.Range("K497").Value = 5000
$ws.Range("L497").Value = 5000
$ws.Range("M497").Value = 5000
$ws.Range("N497").Value = "`$/saco 20 kilos"
$ws.Range("O497").Value = "Región de Ñuble"
$ws.Range("P497").Value = 250
$ws.Range("Q497").Value = 20
$ws.Range("R497").Value = "Hortaliza"
